$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 3
$ws.Range("F7").Value = -6
$ws.Range("F9").Value = -3
$ws.Range("F11").Value = 1
$ws.Range("F20").Value = -2
$ws.Range("F29").Value = 2
$ws.Range("F31").Value = 2
$ws.Range("F34").Value = -1
$ws.Range("F40").Value = -4
$ws.Range("F41").Value = -2
$ws.Range("F43").Value = -3
$ws.Range("F45").Value = 0
$ws.Range("F46").Value = -5
$ws.Range("F51").Value = -4
$ws.Range("F54").Value = -4
$ws.Range("F55").Value = 1
$ws.Range("F56").Value = -4
$ws.Range("F60").Value = -4
$ws.Range("F65").Value = -2
$ws.Range("F67").Value = -7
$ws.Range("F69").Value = -9
$ws.Range("F70").Value = -6
$ws.Range("F71").Value = 6
$ws.Range("F72").Value = -4
$ws.Range("F74").Value = 1
$ws.Range("F76").Value = 4
